$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the "Execute time" values for the two rows of test data.
# These cells currently hold numeric-looking text ("35" / "26") stored as
# shared strings; a leading apostrophe keeps the new values ("34" / "29")
# stored as text too, instead of letting Excel auto-convert them to numbers.
$ws.Range("G2").Value = "'34"
$ws.Range("G3").Value = "'29"
